$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows
$ws.Range("B2").Value = "flower,plants,stem"
$ws.Range("B3").Value = "pot,flower,plants"

# Add new rows
$ws.Range("A4").Value = "cr7.jpg"
$ws.Range("B4").Value = "man, football"

$ws.Range("A5").Value = "C:\Users\Veeraraju_elluru\Downloads\image_tagging_app\image_tagging_app\uploads\l2.jpg"
$ws.Range("B5").Value = "lion, brown, male"

$ws.Range("A6").Value = "C:\Users\Veeraraju_elluru\Downloads\image_tagging_app\image_tagging_app\uploads\m2.jpg"
$ws.Range("B6").Value = "messi, football"

$ws.Range("A7").Value = "a1.jpg"
$ws.Range("B7").Value = "flowers, stem, leaves, blue"
